# Applies the data corrections described by the commit
# "soliucion solicitudes views buscar": several placeholder
# names / ID numbers / addresses / emails inside the acta de
# conciliación are swapped out for different sample values.
#
# Because the same placeholder strings (e.g. "Luz Helena",
# "1234", "Verbenal"...) are reused in more than one place in the
# document -- and not every occurrence changes -- each replacement
# is scoped to the specific paragraph that must change, using that
# paragraph's own Range so the Find cannot leak into neighbouring,
# unrelated paragraphs that happen to contain the same text.
#
# NOTE: this runtime's function parameter binding only works with
# positional arguments (named `-Param value` binding silently fails),
# so Replace-InParagraph takes its arguments positionally.

$d = $word.ActiveDocument

function Replace-InParagraph($ParaIndex, $OldText, $NewText, $WholeWord) {
    $range = $d.Paragraphs.Item($ParaIndex).Range
    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
    #         MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
    #         Format, ReplaceWith, Replace)
    $range.Find.Execute($OldText, $false, $WholeWord, $false, $false, $false, `
                         $true, 1, $false, $NewText, 2) | Out-Null
}

# CONVOCANTE: Luz Helena  ->  Maria Luisa Bogata Rios
Replace-InParagraph 5 "Luz Helena" "Maria Luisa Bogata Rios" $false
#                C.C. No. 1234  ->  123456
Replace-InParagraph 6 "1234" "123456" $true

# CONVOCADO: Kevin Andres Urrego  ->  Maria del Carmen Ruiz C
Replace-InParagraph 8 "Kevin Andres Urrego" "Maria del Carmen Ruiz C" $false
#                C.C. No. 12345  ->  123456789
Replace-InParagraph 9 "12345" "123456789" $true

# Intro paragraph repeats both names / numbers
Replace-InParagraph 11 "Luz Helena" "Maria Luisa Bogata Rios" $false
Replace-InParagraph 11 "1234" "123456" $true
Replace-InParagraph 11 "Kevin Andres Urrego" "Maria del Carmen Ruiz C" $false
Replace-InParagraph 11 "12345" "123456789" $true

# "Estuvieron presentes" - Convocante paragraph
Replace-InParagraph 17 "Luz Helena" "Maria Luisa Bogata Rios" $false
Replace-InParagraph 17 "1234" "123456" $true
Replace-InParagraph 17 "Verbenal" "Las Nieves" $false
Replace-InParagraph 17 "Usaquen" "Santa Fe" $false
Replace-InParagraph 17 "31442169" "12378" $false
Replace-InParagraph 17 "helen@ugc.edu.co" "ro@ugc.edu.co" $false

# Convocada paragraph (Barrio/Localidad stay untouched here)
Replace-InParagraph 19 "Kevin Andres Urrego" "Maria del Carmen Ruiz C" $false
Replace-InParagraph 19 "12345" "123456789" $true
Replace-InParagraph 19 "31247289" "123657657" $false
Replace-InParagraph 19 "kevin@ugc.edu.co" "carmen@ugc.edu.co" $false

# Conciliadora paragraph
Replace-InParagraph 22 "Maria Luisa Bogata Rios" "Jairo Miller Palacio" $false
Replace-InParagraph 22 "123456" "1013689035" $true

# Signature block
Replace-InParagraph 36 "Maria Luisa Bogata Rios" "Jairo Miller Palacio" $false
Replace-InParagraph 37 "123456" "1013689035" $true
Replace-InParagraph 42 "Rafael Perez" "Kevin Andres Urrego" $false
